$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header cell's formatting (bold + border) onto the
# new header cells before overwriting any values.
$ws.Range("B1").Copy()
$ws.Range("C1:D1").PasteSpecial(-4122)  # xlPasteFormats

# Update header row values
$ws.Range("B1").Value = "Exp 1"
$ws.Range("C1").Value = "Exp 2"
$ws.Range("D1").Value = "Exp 3"

# Update data row: label + three numeric results
$ws.Range("A2").Value = "Specificity"
$ws.Range("B2").Value = 0.9692307692307692
$ws.Range("C2").Value = 0.9503546099290781
$ws.Range("D2").Value = 0.9777777777777777
